# Generate Report for Handback
#
# The "c4c46589-1870-45fe-aa05-f073a8d8b456.md" file has now been handed
# back from localization, so update its status from "Ready for handoff"
# to "Handed back: in sync with en-US" on every sheet that tracks it, and
# record the new "Latest Handback DateTime" on the per-locale report
# sheets.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet: collapse the two per-locale status columns ---
$wsOverview = $wb.Worksheets("Overview")
$wsOverview.Range("B3").Value = $handedBack
$wsOverview.Range("C3").Value = $handedBack

# --- zh-cn report sheet ---
$wsZhCn = $wb.Worksheets("zh-cn")
$wsZhCn.Range("B3").Value = $handedBack
$wsZhCn.Range("G3").Value = "2016-02-18 03:41:29"

# --- de-de report sheet ---
$wsDeDe = $wb.Worksheets("de-de")
$wsDeDe.Range("B3").Value = $handedBack
$wsDeDe.Range("G3").Value = "2016-02-18 03:41:51"
